{"js": "// \"MOD. Actualizaci\u00f3n material 2022-2\"\n// Apply the textual edits described by the diff:\n//  1. Typo fix: \"le gustar\u00eda prender de ellos\" -> \"le gustar\u00eda aprender de ellos\"\n//  2. Expand the Google-search sentence with \"Scholar o en Publish or Perish\" (italic)\n//  3. Drop the stray empty (indented / italic) paragraph right before the \"Nota:\" paragraph\n//  4. Insert the \" o un \u201creferente\u201d \" clause into the \"Nota:\" paragraph\n\nconst body = context.document.body;\n\n// 1. \"le gustar\u00eda prender de ellos\" -> \"le gustar\u00eda aprender de ellos\"\n{\n  const results = body.search(\"le gustar\u00eda prender de ellos\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"le gustar\u00eda aprender de ellos\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2. \"...b\u00fasqueda en Google sobre el tema...\" -> \"...b\u00fasqueda en GoogleScholar o en *Publish or Perish* sobre el tema...\"\n{\n  const results = body.search(\"en Google\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const afterGoogle = results.items[0];\n    const scholarRange = afterGoogle.insertText(\"Scholar o en \", Word.InsertLocation.after);\n    await context.sync();\n    const publishRange = scholarRange.insertText(\"Publish or Perish\", Word.InsertLocation.after);\n    publishRange.font.italic = true;\n    await context.sync();\n  }\n}\n\n// 3. Remove the empty paragraph (left indent 1440 twips / 72pt, italic rPr) that used to sit\n//    between \"Despu\u00e9s de este proceso...\" and the \"Nota:\" paragraph.\n{\n  const paras = body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < paras.items.length; i++) {\n    paras.items[i].load(\"text,leftIndent\");\n  }\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    const p = paras.items[i];\n    if (p.text === \"\" && p.leftIndent === 72) {\n      target = p;\n      break;\n    }\n  }\n  if (target) {\n    target.delete();\n    await context.sync();\n  }\n}\n\n// 4. Insert ' o un \"referente\" ' into the \"Nota:\" paragraph, right before \"en ellos, son temas\"\n{\n  const results = body.search(\"en ellos, son temas\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \" o un \\u201creferente\\u201d en ellos, son temas\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# \"MOD. Actualizaci\u00f3n material 2022-2\"\n# Apply the textual edits described by the diff:\n#  1. Typo fix: \"le gustar\u00eda prender de ellos\" -> \"le gustar\u00eda aprender de ellos\"\n#  2. Expand the Google-search sentence with \"Scholar o en Publish or Perish\" (italic)\n#  3. Drop the stray empty (indented / italic) paragraph right before the \"Nota:\" paragraph\n#  4. Insert the \" o un \"referente\" \" clause into the \"Nota:\" paragraph\n\n$d = $word.ActiveDocument\n\n# 1. \"le gustar\u00eda prender de ellos\" -> \"le gustar\u00eda aprender de ellos\"\n$range1 = $d.Content\n$range1.Find.Execute(\"le gustar\u00eda prender de ellos\", $false, $false, $false, $false, $false, $true, 1, $false, \"le gustar\u00eda aprender de ellos\", 2)\n\n# 2. \"...b\u00fasqueda en Google sobre el tema...\" -> \"...b\u00fasqueda en GoogleScholar o en *Publish or Perish* sobre el tema...\"\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\"en Google\")\nif ($found2) {\n    $range2.Collapse(0)                       # wdCollapseEnd\n    $range2.InsertAfter(\"Scholar o en \")\n    $range2.Collapse(0)\n    $range2.InsertAfter(\"Publish or Perish\")\n    $range2.Font.Italic = 1\n}\n\n# 3. Remove the empty paragraph (left indent 1440 twips / 72pt, italic rPr) that used to sit\n#    between \"Despu\u00e9s de este proceso...\" and the \"Nota:\" paragraph.\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text.Length -eq 1 -and $p.Format.LeftIndent -eq 72) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 4. Insert ' o un \"referente\" ' into the \"Nota:\" paragraph, right before \"en ellos, son temas\"\n$range4 = $d.Content\n$found4 = $range4.Find.Execute(\"en ellos, son temas\")\nif ($found4) {\n    $range4.Collapse(1)                       # wdCollapseStart\n    $range4.InsertBefore(\" o un \" + [char]0x201C + \"referente\" + [char]0x201D + \" \")\n}\n"}
